$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

# --- Row 39 / 40 swap: VeChain <-> FraxShare, with updated price/volume ---
Set-TextCell 'B39' 'FraxShare'
Set-TextCell 'C39' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D39' '9.428'
Set-TextCell 'E39' '  +6.50%  '

Set-TextCell 'B40' 'VeChain'
Set-TextCell 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D40' '0.02507'
Set-TextCell 'E40' '  +10.26%  '

# --- Price (D) and Volume 1h (E) updates for remaining rows ---
Set-TextCell 'D2' '25.915.41'
Set-TextCell 'E2' '  +8.98%  '
Set-TextCell 'D3' '1.769.59'
Set-TextCell 'E3' '  +6.90%  '
Set-TextCell 'D4' '1.000'
Set-TextCell 'E4' '  +0.24%  '
Set-TextCell 'D5' '317.32'
Set-TextCell 'E5' '  +3.77%  '
Set-TextCell 'D6' '0.9954'
Set-TextCell 'E6' '  +0.93%  '
Set-TextCell 'D7' '0.3851'
Set-TextCell 'E7' '  +3.74%  '
Set-TextCell 'D8' '0.3658'
Set-TextCell 'E8' '  +6.42%  '
Set-TextCell 'D9' '51.25'
Set-TextCell 'E9' '  +7.49%  '
Set-TextCell 'D10' '1.241'
Set-TextCell 'E10' '  +6.39%  '
Set-TextCell 'D11' '0.07708'
Set-TextCell 'E11' '  +7.31%  '
Set-TextCell 'D12' '0.9968'
Set-TextCell 'E12' '  +0.24%  '
Set-TextCell 'D13' '21.95'
Set-TextCell 'E13' '  +6.99%  '
Set-TextCell 'D14' '6.528'
Set-TextCell 'E14' '  +9.03%  '
Set-TextCell 'D15' '7.156'
Set-TextCell 'E15' '  +6.60%  '
Set-TextCell 'D16' '1.762.35'
Set-TextCell 'E16' '  +6.12%  '
Set-TextCell 'D17' '0.00001172'
Set-TextCell 'E17' '  +7.38%  '
Set-TextCell 'D18' '0.9957'
Set-TextCell 'E18' '  +1.22%  '
Set-TextCell 'D19' '0.06866'
Set-TextCell 'E19' '  +2.07%  '
Set-TextCell 'D20' '87.72'
Set-TextCell 'E20' '  +8.40%  '
Set-TextCell 'D21' '17.79'
Set-TextCell 'E21' '  +8.84%  '
Set-TextCell 'D22' '6.567'
Set-TextCell 'E22' '  +8.38%  '
Set-TextCell 'D23' '12.84'
Set-TextCell 'E23' '  +8.01%  '
Set-TextCell 'D24' '25.853.63'
Set-TextCell 'E24' '  +8.51%  '
Set-TextCell 'D25' '2.432'
Set-TextCell 'E25' '  +4.17%  '
Set-TextCell 'D26' '3.010'
Set-TextCell 'E26' '  +13.09%  '
Set-TextCell 'D27' '20.86'
Set-TextCell 'E27' '  +7.22%  '
Set-TextCell 'D28' '155.46'
Set-TextCell 'E28' '  +2.23%  '
Set-TextCell 'D29' '135.60'
Set-TextCell 'E29' '  +7.15%  '
Set-TextCell 'D30' '1.961.53'
Set-TextCell 'E30' '  +6.40%  '
Set-TextCell 'D31' '1.204'
Set-TextCell 'E31' '  +22.90%  '
Set-TextCell 'D32' '7.282'
Set-TextCell 'E32' '  +19.08%  '
Set-TextCell 'D33' '4.311'
Set-TextCell 'E33' '  +7.09%  '
Set-TextCell 'D34' '14.12'
Set-TextCell 'E34' '  +15.75%  '
Set-TextCell 'D35' '1.817'
Set-TextCell 'E35' '  +6.67%  '
Set-TextCell 'D36' '0.08747'
Set-TextCell 'E36' '  +4.66%  '
Set-TextCell 'D37' '5.718'
Set-TextCell 'E37' '  +8.77%  '
Set-TextCell 'D38' '0.06794'
Set-TextCell 'E38' '  +7.73%  '
Set-TextCell 'D41' '0.2250'
Set-TextCell 'E41' '  +9.62%  '
Set-TextCell 'D42' '1.305'
Set-TextCell 'E42' '  +3.31%  '
Set-TextCell 'D43' '0.6624'
Set-TextCell 'E43' '  +9.71%  '
Set-TextCell 'D44' '14.32'
Set-TextCell 'E44' '  +8.46%  '
Set-TextCell 'D45' '0.9959'
Set-TextCell 'E45' '  +1.33%  '
Set-TextCell 'D46' '0.6406'
Set-TextCell 'E46' '  +8.80%  '
Set-TextCell 'D47' '3.927'
Set-TextCell 'E47' '  +2.77%  '
Set-TextCell 'D48' '2.185'
Set-TextCell 'E48' '  +10.03%  '
Set-TextCell 'D49' '134.13'
Set-TextCell 'E49' '  +5.92%  '
Set-TextCell 'D50' '0.07515'
Set-TextCell 'E50' '  +6.61%  '
Set-TextCell 'D51' '81.31'
Set-TextCell 'E51' '  +7.94%  '
